$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1200
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1200
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1550

$ws.Range("H107").Value = 1955
$ws.Range("I107").Value = 1955
$ws.Range("K107").Value = 1955
$ws.Range("M107").Value = -35

$ws.Range("H112").Value = 998
$ws.Range("J112").Value = 998
$ws.Range("L112").Value = 2994
$ws.Range("N112").Value = -5210

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H137").Value = 994.8333
$ws.Range("I137").Value = 994.8333
$ws.Range("K137").Value = 2984.4999
$ws.Range("M137").Value = -434.4998999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 945.8
$ws.Range("I2").Value = 903.6667
$ws.Range("J2").Value = 1009
$ws.Range("K2").Value = 903.6667
$ws.Range("L2").Value = 1009
$ws.Range("M2").Value = -790.6667
$ws.Range("N2").Value = -1235

$ws.Range("H32").Value = 3620.7778
$ws.Range("I32").Value = 3740.0833
$ws.Range("K32").Value = 3740.0833
$ws.Range("M32").Value = -3453.0833

$ws.Range("H45").Value = 2078.4285
$ws.Range("I45").Value = 2078.4285
$ws.Range("K45").Value = 2078.4285
$ws.Range("M45").Value = -1701.4285

$ws.Range("H102").Value = 2544
$ws.Range("I102").Value = 2544
$ws.Range("K102").Value = 2544
$ws.Range("M102").Value = -922

$ws.Range("H110").Value = 200
$ws.Range("I110").Value = 200
$ws.Range("K110").Value = 200
$ws.Range("M110").Value = 1845

$ws.Range("H116").Value = 945.8
$ws.Range("I116").Value = 903.6667
$ws.Range("J116").Value = 1009
$ws.Range("K116").Value = 903.6667
$ws.Range("L116").Value = 1009
$ws.Range("M116").Value = 1390.3333
$ws.Range("N116").Value = -5597

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 945.8
$ws.Range("I3").Value = 903.6667
$ws.Range("J3").Value = 1009
$ws.Range("K3").Value = 903.6667
$ws.Range("L3").Value = 1009
$ws.Range("M3").Value = -789.6667
$ws.Range("N3").Value = -1237

$ws.Range("H94").Value = 552.3333
$ws.Range("I94").Value = 453.5
$ws.Range("K94").Value = 453.5
$ws.Range("M94").Value = -2.5

$ws.Range("H99").Value = 400
$ws.Range("I99").Value = 400
$ws.Range("K99").Value = 400
$ws.Range("M99").Value = 1098

$ws.Range("H105").Value = 1968.75
$ws.Range("I105").Value = 1968.75
$ws.Range("K105").Value = 1968.75
$ws.Range("M105").Value = -221.75

$ws.Range("H134").Value = 3534.4119
$ws.Range("I134").Value = 3474.0625
$ws.Range("K134").Value = 10422.1875
$ws.Range("M134").Value = -7887.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8605.429
$ws.Range("I16").Value = 8661.25
$ws.Range("J16").Value = 8531
$ws.Range("K16").Value = 8661.25
$ws.Range("L16").Value = 8531
$ws.Range("M16").Value = -8374.25
$ws.Range("N16").Value = -9105

$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1700

$ws.Range("H99").Value = 4999
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H113").Value = 8605.429
$ws.Range("I113").Value = 8661.25
$ws.Range("J113").Value = 8531
$ws.Range("K113").Value = 8661.25
$ws.Range("L113").Value = 8531
$ws.Range("M113").Value = -6491.25
$ws.Range("N113").Value = -12871

$ws.Range("H122").Value = 2025.2
$ws.Range("I122").Value = 1723
$ws.Range("K122").Value = 5169
$ws.Range("M122").Value = -2719

$ws.Range("H126").Value = 4999
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 17973.666
$ws.Range("I132").Value = 41926
$ws.Range("K132").Value = 125778
$ws.Range("M132").Value = -123248

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 430
$ws.Range("I18").Value = 395
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 1185
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -1016
$ws.Range("N18").Value = -1838

$ws.Range("H113").Value = 837.625
$ws.Range("I113").Value = 799.5
$ws.Range("J113").Value = 850.3333
$ws.Range("K113").Value = 2398.5
$ws.Range("L113").Value = 2550.9999
$ws.Range("M113").Value = -228.5
$ws.Range("N113").Value = -6890.9999

$ws.Range("H140").Value = 2032.1428
$ws.Range("I140").Value = 2032.1428
$ws.Range("K140").Value = 6096.428400000001
$ws.Range("M140").Value = -916.4284000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4414.5
$ws.Range("I126").Value = 3497.25
$ws.Range("K126").Value = 10491.75
$ws.Range("M126").Value = -8021.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3863.3333
$ws.Range("I16").Value = 3878
$ws.Range("K16").Value = 3878
$ws.Range("M16").Value = -3708

$ws.Range("H43").Value = 884285.3
$ws.Range("J43").Value = 884285.3
$ws.Range("L43").Value = 884285.3
$ws.Range("N43").Value = -884671.3

$ws.Range("H122").Value = 3264.3333
$ws.Range("I122").Value = 3264.3333
$ws.Range("K122").Value = 9792.999899999999
$ws.Range("M122").Value = -7342.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1165.6666
$ws.Range("I107").Value = 1150
$ws.Range("J107").Value = 1197
$ws.Range("K107").Value = 3450
$ws.Range("L107").Value = 3591
$ws.Range("M107").Value = -1530
$ws.Range("N107").Value = -7431

$ws.Range("H126").Value = 1088
$ws.Range("I126").Value = 443.2
$ws.Range("K126").Value = 1329.6
$ws.Range("M126").Value = 1140.4

$ws.Range("H132").Value = 3216.1304
$ws.Range("I132").Value = 2328.25
$ws.Range("J132").Value = 4184.727
$ws.Range("K132").Value = 6984.75
$ws.Range("L132").Value = 12554.181
$ws.Range("M132").Value = -4454.75
$ws.Range("N132").Value = -17614.181
